# Katpally_LabExam03Grading.xlsx edit:
# Grader entered scores of 10 for E22 and E24 (the "Total Points" column
# for two of the rubric rows in the Table5 grading block, A17:F25).
# The dependent subtotal/total formulas recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E22").Value = 10
$ws.Range("E24").Value = 10

# Move the view/selection to where the grader left off (matches the
# workbook's saved cursor position after entering the scores).
$ws.Range("A11").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 11
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("E24").Select() | Out-Null
